$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
  2  = @(0.011, 6.2, 26.4, 14.3)
  3  = @(-0.055, 15.1, 30.0, 39.5)
  4  = @(-0.044, 19.3, 28.2, 23.7)
  5  = @(0.014, 14.2, 30.5, 4.6)
  6  = @(-0.038, 31.5, 35.0, 25.0)
  7  = @(-0.075, 53.3, 27.5, 7.8)
  8  = @(0.263, 25.1, 21.5, 14.0)
  9  = @(-0.019, 30.2, 24.7, 34.6)
  10 = @(-0.05, 30.2, 27.9, 26.8)
  11 = @(-0.157, 47.3, 22.9, 26.0)
  12 = @(-0.069, 60.8, 17.3, 6.4)
  13 = @(0.189, 34.9, 28.5, 8.9)
  14 = @(-0.04, 46.0, 24.6, 16.2)
  15 = @(-0.13, 60.1, 18.8, 16.2)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("I$row").Value = $vals[0]
    $ws.Range("Y$row").Value = $vals[1]
    $ws.Range("Z$row").Value = $vals[2]
    $ws.Range("AA$row").Value = $vals[3]
}
